$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.593.60'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.65%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.114.25'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.30%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.012'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.91%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '350.61'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +4.88%  '

$ws.Range('E6').Value = '  +0.79%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5248'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.34%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4510'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.33%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '54.41'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +3.07%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.09008'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.03%  '

$ws.Range('E11').Value = '  -0.43%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '24.46'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.18%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '2.118.26'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.97%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.822'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.08%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '8.029'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.72%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '101.22'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +4.86%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001170'
$ws.Range('D17').Style = 'Normal'

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '1.011'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.76%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06711'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.30%  '

$ws.Range('E20').Value = '  +0.73%  '

$ws.Range('E21').Value = '  +0.76%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.297'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.12%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '30.673.02'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.70%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '12.84'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +3.27%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.396'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.44%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.375.38'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.13%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '22.40'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.31%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '165.09'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.89%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.538'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.37%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '136.13'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +2.57%  '

$ws.Range('E31').Value = '  -3.93%  '

$ws.Range('E32').Value = '  +0.44%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.645'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -5.30%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.371'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.95%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.019'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +2.54%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '10.37'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.70%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.934'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +6.76%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02648'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +2.46%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.06841'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.11%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.2319'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.95%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '12.53'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.09%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.6876'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.39%  '

$ws.Range('E43').Value = '  +1.89%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '14.68'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +4.49%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.332'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.06%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.6448'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.06%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.758'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.72%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.00000000357'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.29%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.252'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.41%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.07290'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.37%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '82.48'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.41%  '

